# Edit "TextBox 4" on slide 2:
#  - shrink the shape's height
#  - reduce the font size of both paragraphs (48pt -> 44pt)
#  - add a new "YouTube" entry to the second line of text

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item("TextBox 4")

$tf = $shp.TextFrame
$tr = $tf.TextRange

# --- Paragraph 1: "Wiki \u2022 Forums \u2022 GitHub" -> just a font-size change
$para1 = $tr.Paragraphs(1)
$para1.Font.Size = 44

# --- Paragraph 2: "Quora \u2022 RosettaCode" -> "\tYouTube \u2022 Quora \u2022 RosettaCode"
$para2 = $tr.Paragraphs(2)

# Replace the leading "Quora " text with the new "\tYouTube \u2022 Quora " text;
# this keeps the following "\u2022 RosettaCode" runs intact (same run-split
# behaviour PowerPoint itself performs on a partial retype).
$lead = $tr.Characters($para2.Start, 6)
$lead.Text = [char]9 + "YouTube " + [char]0x2022 + " Quora "

# Bring the whole (now 3-run) second paragraph down to the new font size
$para2 = $tr.Paragraphs(2)
$para2.Font.Size = 44

# The textbox has <a:spAutoFit/> - shrinking the font means the rendered
# (autofit) height shrinks too. Re-create that laid-out extent precisely
# (cx/left/top are untouched - only the height changes).
$targetHeightEmu = 1446550
$shp.Height = [math]::Round($targetHeightEmu / 914400 * 72, 4)
